$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set D22 first so its shared string ("Implementation of Round-Robin...")
# becomes shared-string index 21, matching the order new strings are
# appended in the target workbook.
$ws.Range("D22").Value = "Implementation of Round-Robin. GNU FDL added to manual"
$ws.Range("D21").Value = "Manual set up from LaTeX template"

# Row 21: Manual set up from LaTeX template
$ws.Range("A21").Value = 41190
$ws.Range("B21").Value = 1

# Row 22: Implementation of Round-Robin. GNU FDL added to manual
$ws.Range("A22").Value = 41192
$ws.Range("B22").Value = 2

# Update the active selection to match the post-edit state
$ws.Range("A23").Select()

$wb.Save()
